# Edit script: "Adionamento dos artefatos 16-17-19"
#
# The document has two tables describing stakeholders of the RoadOn
# software. Two bullet items in the "Empresas de turismo" row of the
# first table need their wording updated:
#
#   1) "Promover programas de fidelidade;"
#        -> "Gerenciamento de clientes;"
#
#   2) "...desde o aluguel de ônibus ou vans, até o fim da viagem -,..."
#        -> "...desde o início da viagem até o final da mesma -,..."
#
$d = $word.ActiveDocument

# 1) Update the bullet about organizing trips: replace the travel-range
#    description in the middle of the sentence.
$d.Content.Find.Execute(
    "desde o aluguel de ônibus ou vans, até o fim da viagem -,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "desde o início da viagem até o final da mesma -,",
    2)

# 2) Replace the "loyalty programs" bullet with "client management".
$d.Content.Find.Execute(
    "Promover programas de fidelidade;",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Gerenciamento de clientes;",
    2)
